$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "27.764.47"
$ws.Cells.Item(2, 5).Value = "  +0.12%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.622.09"
$ws.Cells.Item(3, 5).Value = "  -0.16%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.998"
$ws.Cells.Item(4, 5).Value = "  +0.55%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "210.71"
$ws.Cells.Item(5, 5).Value = "  +0.08%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.518"
$ws.Cells.Item(6, 5).Value = "  -0.54%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.997"
$ws.Cells.Item(7, 5).Value = "  +0.52%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "23.12"
$ws.Cells.Item(8, 5).Value = "  +0.03%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.255"
$ws.Cells.Item(9, 5).Value = "  -0.65%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0607"
$ws.Cells.Item(10, 5).Value = "  -0.59%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0878"
$ws.Cells.Item(11, 5).Value = "  +0.05%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.851.02"
$ws.Cells.Item(12, 5).Value = "  -0.18%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.661.65"
$ws.Cells.Item(13, 5).Value = "  +2.68%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.00"
$ws.Cells.Item(14, 5).Value = "  -0.77%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.556"
$ws.Cells.Item(15, 5).Value = "  -0.65%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "64.92"
$ws.Cells.Item(16, 5).Value = "  -0.41%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "27.775.90"
$ws.Cells.Item(17, 5).Value = "  +0.17%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "226.62"
$ws.Cells.Item(18, 5).Value = "  -1.95%  "
$ws.Cells.Item(19, 2).Value = "ShibaInu"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.0₃0717"
$ws.Cells.Item(19, 5).Value = "  -0.60%  "
$ws.Cells.Item(20, 2).Value = "Chainlink"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.58"
$ws.Cells.Item(20, 5).Value = "  +1.45%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.997"
$ws.Cells.Item(21, 5).Value = "  +0.57%  "
$ws.Cells.Item(22, 5).Value = "  -0.46%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.92"
$ws.Cells.Item(23, 5).Value = "  -3.87%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.04"
$ws.Cells.Item(24, 5).Value = "  +0.02%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "155.16"
$ws.Cells.Item(25, 5).Value = "  +1.11%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "6.91"
$ws.Cells.Item(26, 5).Value = "  +0.42%  "
$ws.Cells.Item(27, 5).Value = "  -0.12%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "15.41"
$ws.Cells.Item(28, 5).Value = "  -0.95%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.998"
$ws.Cells.Item(29, 5).Value = "  +0.46%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.18"
$ws.Cells.Item(30, 5).Value = "  +0.29%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.0479"
$ws.Cells.Item(31, 5).Value = "  +0.17%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.38"
$ws.Cells.Item(32, 5).Value = "  -0.09%  "
$ws.Cells.Item(33, 5).Value = "  +0.57%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.402.98"
$ws.Cells.Item(34, 5).Value = "  +0.58%  "
$ws.Cells.Item(35, 5).Value = "  +2.79%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.00"
$ws.Cells.Item(36, 5).Value = "  -0.23%  "
$ws.Cells.Item(37, 5).Value = "  -0.28%  "
$ws.Cells.Item(38, 5).Value = "  +0.19%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.554"
$ws.Cells.Item(39, 5).Value = "  -0.27%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.842"
$ws.Cells.Item(40, 5).Value = "  -2.40%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.997"
$ws.Cells.Item(41, 5).Value = "  +0.52%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.999"
$ws.Cells.Item(42, 5).Value = "  -1.66%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.83"
$ws.Cells.Item(43, 5).Value = "  +0.94%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "65.55"
$ws.Cells.Item(44, 5).Value = "  -0.99%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "5.41"
$ws.Cells.Item(45, 5).Value = "  -0.20%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.761.04"
$ws.Cells.Item(46, 5).Value = "  -0.24%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.11"
$ws.Cells.Item(47, 5).Value = "  -3.34%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "87.91"
$ws.Cells.Item(48, 5).Value = "  +0.22%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.101"
$ws.Cells.Item(49, 5).Value = "  +1.17%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0502"
$ws.Cells.Item(50, 5).Value = "  -0.39%  "
$ws.Cells.Item(51, 5).Value = "  +1.22%  "

Write-Host "Done"